$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B12").Value = 22
$ws.Range("B13").Value = 770000
$ws.Range("B14").Value = 2357142.857142857
$ws.Range("B33").Value = 1797142.857142857
$ws.Range("B35").Value = 1797142.857142857

$wb.Save()
